# Revert "Revert "Add Asmara"" — re-insert the "hh" / "Asmara" row into the
# Owners sheet, in alphabetical position right after the "hc" / "Mogadishu"
# row (row 44), pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 45 (everything from old row 45 onward
# shifts down to 46+). Excel's Insert() carries the formatting of the row
# above down into the new row, which already matches what we need here
# (style 1 for A:C, style 2 / merged-looking for D:E).
$ws.Rows("45:45").Insert()

# Populate the new "hh" / Asmara entry.
$ws.Range("A45").Value = "hh"
$ws.Range("B45").Value = "Asmara"
$ws.Range("C45").Value = "VATGlasses"
$ws.Range("D45").Value = "vACC Rejected Offer To Participate"

# D:E are merged on every "VATGlasses" row in this table.
$ws.Range("D45:E45").Merge()

# Restore the selection/scroll position to match the post-edit workbook
# (user had scrolled to and selected the new row's D:E cell).
$ws.Range("D45:E45").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
